$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 786, shifting existing rows 786-840 down to 787-841
$ws.Rows("786").Insert()

# Populate the newly inserted row 786 with the new record
$ws.Range("A786").Value = 5
$ws.Range("B786").Value = "Macroferia Regional de Talca"
$ws.Range("C786").Value = "Maule"
$ws.Range("D786").Value = 45021
$ws.Range("E786").Value = 7
$ws.Range("F786").Value = 100112004
$ws.Range("G786").Value = "Cebolla"
$ws.Range("H786").Value = "Sin especificar"
$ws.Range("I786").Value = "1a (cosecha)"
$ws.Range("J786").Value = 2500
$ws.Range("K786").Value = 10000
$ws.Range("L786").Value = 10000
$ws.Range("M786").Value = 10000
$ws.Range("N786").Value = "$/malla 25 kilos"
$ws.Range("O786").Value = "Región del Maule"
$ws.Range("P786").Value = 400
$ws.Range("Q786").Value = 25
$ws.Range("R786").Value = "Hortaliza"
